# Applies the cryptos-list price/volume/coin updates described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "69.118.89"
$ws.Range("E2").Value = "  +0.51%  "
# Row 3
$ws.Range("D3").Value = "3.752.98"
$ws.Range("E3").Value = "  +0.65%  "
# Row 4
$ws.Range("E4").Value = "  +0.17%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.76"
$ws.Range("E5").Value = "  +0.01%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.26"
$ws.Range("E6").Value = "  -1.99%  "
# Row 7
$ws.Range("D7").Value = "3.753.81"
$ws.Range("E7").Value = "  +0.66%  "
# Row 8
$ws.Range("E8").Value = "  +0.00%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.540"
$ws.Range("E9").Value = "  +0.92%  "
# Row 10
$ws.Range("E10").Value = "  +4.95%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.41"
$ws.Range("E11").Value = "  +1.10%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.460"
$ws.Range("E12").Value = "  -0.30%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.83"
$ws.Range("E13").Value = "  -0.85%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000248"
$ws.Range("E14").Value = "  +1.13%  "
# Row 15
$ws.Range("D15").Value = "4.381.21"
$ws.Range("E15").Value = "  +0.65%  "
# Row 16
$ws.Range("D16").Value = "3.754.83"
$ws.Range("E16").Value = "  +0.85%  "
# Row 17
$ws.Range("D17").Value = "69.239.41"
$ws.Range("E17").Value = "  +0.78%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.44"
$ws.Range("E18").Value = "  +2.15%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.43"
$ws.Range("E19").Value = "  +1.56%  "
# Row 20
$ws.Range("E20").Value = "  -1.78%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.14"
$ws.Range("E21").Value = "  +6.32%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "493.37"
$ws.Range("E22").Value = "  -0.58%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.727"
$ws.Range("E23").Value = "  +0.29%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000149"
$ws.Range("E24").Value = "  +3.60%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.88"
$ws.Range("E25").Value = "  -0.33%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.28"
$ws.Range("E26").Value = "  -1.57%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.36"
$ws.Range("E27").Value = "  -0.95%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.04"
$ws.Range("E28").Value = "  -0.98%  "
# Row 29
$ws.Range("E29").Value = "  -0.10%  "
# Row 30
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.97"
$ws.Range("E30").Value = "  +0.17%  "
# Row 31
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.23"
$ws.Range("E31").Value = "  +3.22%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.47"
$ws.Range("E32").Value = "  -4.29%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.69"
$ws.Range("E33").Value = "  +0.27%  "
# Row 34
$ws.Range("D34").Value = "3.900.64"
$ws.Range("E34").Value = "  +0.69%  "
# Row 35
$ws.Range("D35").Value = "3.686.35"
$ws.Range("E35").Value = "  +0.62%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.108"
$ws.Range("E36").Value = "  -0.20%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.97"
$ws.Range("E37").Value = "  +2.15%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.139"
$ws.Range("E38").Value = "  +6.03%  "
# Row 39
$ws.Range("E39").Value = "  -0.33%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.20"
$ws.Range("E40").Value = "  +10.83%  "
# Row 41
$ws.Range("E41").Value = "  +0.07%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.326"
$ws.Range("E42").Value = "  +0.18%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.00"
$ws.Range("E43").Value = "  +0.55%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "48.63"
$ws.Range("E44").Value = "  -0.47%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "426.45"
$ws.Range("E45").Value = "  -2.91%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.44"
$ws.Range("E46").Value = "  -0.56%  "
# Row 47
$ws.Range("E47").Value = "  -0.01%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.11"
$ws.Range("E48").Value = "  -1.24%  "
# Row 49
$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.31"
$ws.Range("E49").Value = "  +7.61%  "
# Row 50
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "141.11"
$ws.Range("E50").Value = "  -0.01%  "
# Row 51
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "2.790.90"
$ws.Range("E51").Value = "  +1.10%  "
